$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A163").Value = "DEVOPS LEAD WITH Google Cloud Platform, KUBERNETES, TERRAFORM"
$ws.Range("B163").Value = "https://www.dice.com/job-detail/fde541bd-0775-47ec-90eb-44325d8159ef?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=3&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C163").Value = "San Jose, California"
$ws.Range("D163").Value = "Third Party, Contract"
$ws.Range("E163").Value = "Depends on Experience"
$ws.Range("F163").Value = "Emergere Technologies"

$ws.Range("A164").Value = "Java Full Stack Backend Developer"
$ws.Range("B164").Value = "https://www.dice.com/job-detail/23c947c6-5f82-4e36-922c-933d8a6f5d3c?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=5&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
$ws.Range("C164").Value = "San Leandro, California"
$ws.Range("D164").Value = "Third Party, Contract"
$ws.Range("E164").Value = "Depends on Experience"
$ws.Range("F164").Value = "Source Code Technologies LLC"
